# Update column G ("K" - strikeouts) values on Sheet1 with regenerated
# data (using K instead of Strike#). Only column G values change; all
# other columns / formatting remain untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newK = @{
    2  = 6
    3  = 3
    4  = 14
    5  = 6
    6  = 9
    7  = 6
    8  = 4
    9  = 4
    10 = 3
    11 = 9
    12 = 8
    13 = 16
    14 = 10
    15 = 8
    16 = 4
    17 = 7
    18 = 11
    19 = 6
    20 = 1
    21 = 10
    22 = 5
    23 = 8
    24 = 2
    25 = 7
    26 = 3
    27 = 0
    28 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
